$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("joghurt, natúr joghurt" / 60) is removed from the foods table.
# Deleting the whole row shifts every row below it up by one, which also
# keeps the A2:A.. list sorted alphabetically (matches the rest of the
# table - nothing else needs to be re-sorted).
$ws.Rows.Item(16).Delete()

# A brand-new row is appended at the bottom of the table for "saláta" (300 kcal/100g).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$ws.Cells.Item($lastRow, 1).Value = "saláta"
$ws.Cells.Item($lastRow, 2).Value = 300

# Restore/update the view: scroll so row 13 is at the top and select B50
# (the cell just below the newly appended row), matching the saved workbook view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B50").Select()
